# Apply "Finished updates to table 2" - refreshed SIBER stats values,
# normalized the Probability column (H) to the same 2-decimal number
# format as the rest of the table (was a custom 4-decimal "0.0000"),
# and cleared the explicit right-alignment override on the data cells
# (columns D:I) since the column-level default already right-aligns them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated statistics values (columns D:I, rows 2-9) ---------------

# Row 2 - a) Quaternary Sylvilagus
$ws.Range("D2").Value = 2.7829088855611102
$ws.Range("E2").Value = 2.85248160770014
$ws.Range("F2").Value = 2.0155688298500301
$ws.Range("G2").Value = 3.7247299855407299
$ws.Range("H2").Value = 0.98250000000000004
$ws.Range("I2").Value = 0.80806077549286803

# Row 3 - Quaternary Otospermophilus
$ws.Range("D3").Value = 4.7104154385561499
$ws.Range("E3").Value = 4.8988320560984002
$ws.Range("F3").Value = 3.0840138438368498
$ws.Range("G3").Value = 6.8220137036014998
$ws.Range("I3").Value = 0.47051592575588502

# Row 4 - b) Pleistocene Sylvilagus
$ws.Range("D4").Value = 1.9092959765366999
$ws.Range("E4").Value = 1.95954060749819
$ws.Range("F4").Value = 1.3672246309695999
$ws.Range("G4").Value = 2.5997985511004602
$ws.Range("H4").Value = 0.97824999999999995
$ws.Range("I4").Value = 0.73495129108706903

# Row 5 - Pleistocene Otospermophilus
$ws.Range("D5").Value = 3.21963437403444
$ws.Range("E5").Value = 3.3659813910360099
$ws.Range("F5").Value = 2.0834212831486201
$ws.Range("G5").Value = 4.7669781319710696
$ws.Range("I5").Value = 0.42785943595938702

# Row 6 - c) All Taxa Pleistocene
$ws.Range("D6").Value = 3.1293811346673102
$ws.Range("E6").Value = 3.1798550239361401
$ws.Range("F6").Value = 2.4093022881636501
$ws.Range("G6").Value = 3.95788454649094
$ws.Range("H6").Value = 0.99124999999999996
$ws.Range("I6").Value = 0.57429633114927303

# Row 7 - All Taxa Holocene
$ws.Range("F7").Value = 2.34422685486456
$ws.Range("G7").Value = 19.3132212421236
$ws.Range("I7").Value = 0.19641445360831

# Row 8 - d) All Taxa Pleistocene
$ws.Range("D8").Value = 3.1293811346673102
$ws.Range("E8").Value = 3.1798550239361401
$ws.Range("F8").Value = 2.4093022881636501
$ws.Range("G8").Value = 3.95788454649094
$ws.Range("H8").Value = 0.95825000000000005

# Row 9 - All Taxa Quaternary
$ws.Range("D9").Value = 4.24066348004588
$ws.Range("E9").Value = 4.3039569648226896
$ws.Range("F9").Value = 3.3089007975008502
$ws.Range("G9").Value = 5.34136658569318

# --- Number format: Probability column no longer uses the custom ----
# --- "0.0000" format; it now matches the other numeric columns ------
$ws.Range("H2:H9").NumberFormat = "0.00"

# --- Clear the explicit right-alignment that was baked into the -----
# --- per-cell style of the data columns (column default still -------
# --- right-aligns numbers, so this is purely a style cleanup) -------
$ws.Range("D2:I9").HorizontalAlignment = 1

# --- Leave the cursor where the author left it when they saved ------
$ws.Range("C11").Select()
